$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row labels first (column A, rows 2-5) -- matches shared-string order 0..3
$ws.Cells.Item(2, 1).Value = "Total Points"
$ws.Cells.Item(3, 1).Value = "Games Played"
$ws.Cells.Item(4, 1).Value = "Wins"
$ws.Cells.Item(5, 1).Value = "Losses"

# Header row (row 1) -- matches shared-string order 4..8
$ws.Cells.Item(1, 1).Value = "Players"
$ws.Cells.Item(1, 2).Value = "Nick"
$ws.Cells.Item(1, 3).Value = "Matt"
$ws.Cells.Item(1, 4).Value = "Jasper"
$ws.Cells.Item(1, 5).Value = "Nolan"

# Numeric data
$data = @(
    @(150, 100, 190, 135),
    @(10, 10, 10, 10),
    @(4, 5, 8, 1),
    @(6, 5, 2, 9)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}

$ws.PageSetup.Orientation = 1

$ws.Range("A5").Select()
